$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.302.56"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "3.921.45"
$ws.Range("E3").Value = "  -1.38%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'487.21"
$ws.Range("E5").Value = "  +0.45%  "

$ws.Range("D6").Value = "'146.40"
$ws.Range("E6").Value = "  -2.21%  "

$ws.Range("E7").Value = "  -1.01%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("D9").Value = "'0.734"
$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("E10").Value = "  -1.80%  "

$ws.Range("D11").Value = "'0.0000350"
$ws.Range("E11").Value = "  -5.18%  "

$ws.Range("D12").Value = "'43.21"
$ws.Range("E12").Value = "  -1.41%  "

$ws.Range("D13").Value = "'10.72"
$ws.Range("E13").Value = "  +1.94%  "

$ws.Range("D14").Value = "4.540.25"
$ws.Range("E14").Value = "  -1.34%  "

$ws.Range("D15").Value = "3.917.56"
$ws.Range("E15").Value = "  -1.44%  "

$ws.Range("D16").Value = "'14.30"
$ws.Range("E16").Value = "  -4.09%  "

$ws.Range("E17").Value = "  -0.79%  "

$ws.Range("D18").Value = "'20.12"
$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("E19").Value = "  -0.96%  "

$ws.Range("D20").Value = "68.350.72"
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").Value = "'432.12"
$ws.Range("E21").Value = "  -1.45%  "

$ws.Range("D22").Value = "'3.51"
$ws.Range("E22").Value = "  +3.18%  "

$ws.Range("D23").Value = "'15.17"
$ws.Range("E23").Value = "  +4.96%  "

$ws.Range("D24").Value = "'88.57"
$ws.Range("E24").Value = "  +0.05%  "

$ws.Range("D25").Value = "'11.78"
$ws.Range("E25").Value = "  +22.00%  "

$ws.Range("D26").Value = "'3.72"
$ws.Range("E26").Value = "  +1.77%  "

$ws.Range("D27").Value = "'11.21"
$ws.Range("E27").Value = "  +10.47%  "

$ws.Range("D28").Value = "'37.89"
$ws.Range("E28").Value = "  -2.70%  "

$ws.Range("E29").Value = "  -1.62%  "

$ws.Range("D30").Value = "'718.78"
$ws.Range("E30").Value = "  -1.86%  "

$ws.Range("D31").Value = "'13.76"
$ws.Range("E31").Value = "  +3.11%  "

$ws.Range("E32").Value = "  +0.84%  "

$ws.Range("D33").Value = "'2.92"
$ws.Range("E33").Value = "  +2.89%  "

$ws.Range("D34").Value = "0.0₃0917"
$ws.Range("E34").Value = "  +3.37%  "

$ws.Range("D35").Value = "'6.17"
$ws.Range("E35").Value = "  +14.50%  "

$ws.Range("D36").Value = "'41.78"
$ws.Range("E36").Value = "  -0.74%  "

$ws.Range("D37").Value = "'61.02"
$ws.Range("E37").Value = "  +0.62%  "

$ws.Range("D38").Value = "'0.398"
$ws.Range("E38").Value = "  +18.39%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.146"
$ws.Range("E39").Value = "  -3.91%  "

$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "'3.01"
$ws.Range("E41").Value = "  +18.24%  "

$ws.Range("D42").Value = "'0.0493"
$ws.Range("E42").Value = "  +3.73%  "

$ws.Range("E43").Value = "  +2.90%  "

$ws.Range("E44").Value = "  +4.71%  "

$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("D46").Value = "'3.35"
$ws.Range("E46").Value = "  +2.92%  "

$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").Value = "'3.43"
$ws.Range("E48").Value = "  -0.49%  "

$ws.Range("E49").Value = "  -5.52%  "

$ws.Range("D50").Value = "'144.96"
$ws.Range("E50").Value = "  -2.59%  "

$ws.Range("E51").Value = "  +26.86%  "

